# Update the "NewUsers" sheet test data: rename Shwetha_55 -> Shwetha_56
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewUsers")

$ws.Range("A2").Value = "Shwetha_56"
$ws.Range("B2").Value = "Shwetha_56@gmail.om"

# Move the active selection on this sheet to E14 (matches the saved view state)
$ws.Activate()
$ws.Range("E14").Select()
